$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for every data row (2..267).
$ws.Range("C2:C267").Value = 45182
